$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:P2")
$rng.Style = "Normal"

$ws.Range("B2").Value = 12343555

$ws.Range("B6").Select()
